# Adds Faculty details to the ECE semester-7 timetable workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Course_Information: insert a "Faculty" column (new col E) between
#    "Term Type"/"Basket" (D) and "Display Format" (which shifts E -> F).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Course_Information")
$ws.Columns("E:E").Insert()
$ws.Columns("E:E").ColumnWidth = 24.17

$ws.Range("E4").Value = "Faculty"
$ws.Range("E5").Value = "Pankaj Kumar"

$ws.Range("E8").Value = "Faculty"
$ws.Range("E9").Value = "Sandesh P"
$ws.Range("E10").Value = "Shirshendu Layek"
$ws.Range("E11").Value = "Krishnendu"
$ws.Range("E12").Value = "Abdul Wahid"
$ws.Range("E13").Value = "Malay Kumar"
$ws.Range("E14").Value = "Sandesh Phalke"
$ws.Range("E15").Value = "Anushree"
$ws.Range("E16").Value = "Girish G N"
$ws.Range("E17").Value = "Rajendra Hegadi"
$ws.Range("E18").Value = "Sunil Saumya"
$ws.Range("E19").Value = "Dibyajyothi"
$ws.Range("E20").Value = "Chinmayananda A"
$ws.Range("E21").Value = "Jagadish D N"
$ws.Range("E22").Value = "Rajesh Kumar"
$ws.Range("E23").Value = "Anand B"
$ws.Range("E24").Value = "Aswath Babu"

# ---------------------------------------------------------------------------
# 2. Regular_Timetable / PreMid_Timetable / PostMid_Timetable: the EC462
#    "Generative AI" mini-project's room assignment changed so several
#    elective rows' room codes in columns D/E were re-pointed to share a
#    single room per course instead of separate rooms.
# ---------------------------------------------------------------------------
$timetableSheets = @("Regular_Timetable", "PreMid_Timetable", "PostMid_Timetable")
foreach ($sheetName in $timetableSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("D20").Value = "Mon 09:00-10:30 [C101], Wed 13:00-14:30 [C101]"
    $ws.Range("E20").Value = "Tue 14:30-15:30 [C101]"

    $ws.Range("D21").Value = "Mon 09:00-10:30 [C102], Wed 13:00-14:30 [C102]"
    $ws.Range("E21").Value = "Tue 14:30-15:30 [C102]"

    $ws.Range("D22").Value = "Mon 09:00-10:30 [C104], Wed 13:00-14:30 [C104]"
    $ws.Range("E22").Value = "Tue 14:30-15:30 [C104]"

    $ws.Range("D23").Value = "Tue 09:00-10:30 [C101], Thu 13:00-14:30 [C101]"
    $ws.Range("E23").Value = "Wed 14:30-15:30 [C101]"

    $ws.Range("D24").Value = "Tue 09:00-10:30 [C102], Thu 13:00-14:30 [C102]"
    $ws.Range("E24").Value = "Wed 14:30-15:30 [C102]"

    $ws.Range("D25").Value = "Tue 09:00-10:30 [C104], Thu 13:00-14:30 [C104]"
    $ws.Range("E25").Value = "Wed 14:30-15:30 [C104]"

    $ws.Range("D26").Value = "Tue 09:00-10:30 [C202], Thu 13:00-14:30 [C202]"
    $ws.Range("E26").Value = "Wed 14:30-15:30 [C202]"

    $ws.Range("E27").Value = "Thu 14:30-15:30 [C101]"
    $ws.Range("E28").Value = "Thu 14:30-15:30 [C102]"
    $ws.Range("E29").Value = "Thu 14:30-15:30 [C104]"
    $ws.Range("E30").Value = "Thu 14:30-15:30 [C202]"
    $ws.Range("E31").Value = "Thu 14:30-15:30 [C203]"
}

# ---------------------------------------------------------------------------
# 3. Section_A: the Mini Project (EC498) room moved from C001 to C204.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Section_A")
$ws.Range("B16").Value = "Mini Project [C204]"
$ws.Range("C16").Value = "0-0-0-8-2 [C204]"
$ws.Range("D16").Value = "Full Sem [C204]"
$ws.Range("E16").Value = "0/0 [C204]"
$ws.Range("F16").Value = "0/0 [C204]"

# ---------------------------------------------------------------------------
# 4. Classroom_Utilization: the logged weekly/daily hours move from room
#    C001 (row 2) to room C204 (row 16).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Classroom_Utilization")
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("D16").Value = 7.5
$ws.Range("E16").Value = 1.5

# ---------------------------------------------------------------------------
# 5. Classroom_Allocation: EC498 Mini Project moves from C001 to C204
#    (and gets reclassified from a "large classroom"/120 to a regular
#    "classroom"/96 with a TV); several electives' TV-equipped rooms were
#    swapped to Projector-equipped rooms, and rooms were consolidated so
#    each course only uses a single room across its sessions.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Classroom_Allocation")

# "96" must stay a text value (it mirrors the neighbouring Capacity column
# which is stored as text) - writing the literal numeric-looking string via
# .Value would get auto-coerced to a number, so stage it as text in a
# scratch cell and paste only the *value* across, which keeps it a string
# without carrying along the scratch cell's text-number-format style.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "96"
$scratch.Copy()
foreach ($r in @(3, 6, 9, 12, 15, 35)) {
    $ws.Range("H$r").PasteSpecial(-4163)
}
$scratch.Clear()

foreach ($r in @(3, 6, 9, 12, 15)) {
    $ws.Range("G$r").Value = "classroom"
    $ws.Range("I$r").Value = "TV"
    $ws.Range("M$r").Value = "C204"
}

$ws.Range("M31").Value = "C104"

$ws.Range("I32").Value = "Projector"
$ws.Range("M32").Value = "C202"

$ws.Range("I33").Value = "Projector"
$ws.Range("M33").Value = "C101"

$ws.Range("I34").Value = "Projector"
$ws.Range("M34").Value = "C102"

$ws.Range("G35").Value = "classroom"
$ws.Range("I35").Value = "Projector"
$ws.Range("M35").Value = "C104"

$ws.Range("I47").Value = "Projector"
$ws.Range("M47").Value = "C202"

$ws.Range("M48").Value = "C203"

$ws.Range("M49").Value = "C101"

$ws.Range("I50").Value = "Projector"
$ws.Range("M50").Value = "C102"

$ws.Range("I51").Value = "Projector"
$ws.Range("M51").Value = "C104"

$ws.Range("I52").Value = "Projector"
$ws.Range("M52").Value = "C101"

$ws.Range("I53").Value = "Projector"
$ws.Range("M53").Value = "C102"

$ws.Range("I54").Value = "Projector"
$ws.Range("M54").Value = "C104"

$ws.Range("M55").Value = "C202"

$ws.Range("I60").Value = "Projector"
$ws.Range("M60").Value = "C101"

$ws.Range("I61").Value = "Projector"
$ws.Range("M61").Value = "C102"

$ws.Range("I62").Value = "Projector"
$ws.Range("M62").Value = "C104"

$ws.Range("M63").Value = "C202"

$ws.Range("I64").Value = "Projector"
$ws.Range("M64").Value = "C101"

$ws.Range("M65").Value = "C102"

$ws.Range("I66").Value = "Projector"
$ws.Range("M66").Value = "C104"

$ws.Range("I67").Value = "Projector"
$ws.Range("M67").Value = "C202"

$ws.Range("M68").Value = "C203"

# ---------------------------------------------------------------------------
# 6. Basket_Course_Allocations: each elective course's candidate room list
#    is trimmed down to the room(s) it actually ended up using.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Basket_Course_Allocations")
$ws.Range("C2").Value = "C004, C101"
$ws.Range("C3").Value = "C102"
$ws.Range("C4").Value = "C104"
$ws.Range("C5").Value = "C004, C101"
$ws.Range("C6").Value = "C102"
$ws.Range("C7").Value = "C104"
$ws.Range("C8").Value = "C202"
$ws.Range("C9").Value = "C203"
$ws.Range("C10").Value = "C004, C101"
$ws.Range("C11").Value = "C102"
$ws.Range("C12").Value = "C104"
$ws.Range("C13").Value = "C202"

Write-Host "Applied Faculty details + room-consolidation updates"
